# Fruta / hortaliza, semanal
# Insert a new weekly record at row 94 (pushing the existing rows 94-107
# down to 95-108) on the "Hortaliza, Vega Modelo de Temuco - Espárragos"
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 94:107 down by one to make room for the new record; Excel
# carries the row-94 formatting (the date-formatted column D style) onto
# the freshly inserted blank row, matching the existing rows.
$ws.Rows("94:94").Insert()

# Populate the new row 94 with the latest weekly price data.
$ws.Range("A94").Value = 10
$ws.Range("B94").Value = "Vega Modelo de Temuco"
$ws.Range("C94").Value = "La Araucanía"
$ws.Range("D94").Value = 45223
$ws.Range("E94").Value = 9
$ws.Range("F94").Value = 300000000
$ws.Range("G94").Value = "Espárragos"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 900
$ws.Range("K94").Value = 1500
$ws.Range("L94").Value = 1600
$ws.Range("M94").Value = 1556
$ws.Range("N94").Value = "$/kilo"
$ws.Range("O94").Value = "Región del Maule"
$ws.Range("P94").Value = 1556
$ws.Range("Q94").Value = 1
$ws.Range("R94").Value = "Hortaliza"
